# Rename sheet "Property1" -> "DataNode" to match the unified
# DataNode/DataTable/Entity naming convention referenced in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "DataNode"

# Move the selection on the frozen (bottom-left) pane to F25, matching the
# saved view state captured in the workbook.
$ws.Range("F25").Select() | Out-Null
